$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.236299999999996
$ws.Range("D8").Value = -8.745399999999991
$ws.Range("A12").Value = -22.88230000000001
$ws.Range("D12").Value = -8.211
$ws.Range("D14").Value = -8.708499999999997
$ws.Range("D22").Value = -7.733800000000002
